# Add "Direction" and "Asset Name" columns to the EntryPoint table, and add a
# new "Vectors" worksheet/table at the end of the workbook.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------------
# 1. EntryPoint sheet: insert "Direction" (after ID) and "Asset Name"
#    (after Description) columns into the Table2 listobject.
# -----------------------------------------------------------------------
$epSheet = $wb.Worksheets.Item("EntryPoint")
$epTable = $epSheet.ListObjects.Item(1)

# Insert the two new columns as plain worksheet columns first (this keeps
# the widths/bestFit of the untouched columns intact), then resize the
# table to cover them.
$epSheet.Columns.Item(2).Insert()   # new column B, "Direction"
$epSheet.Columns.Item(4).Insert()   # new column D, "Asset Name"
$epSheet.Columns.Item(4).ClearFormats()

$epSheet.Range("B1").Value = "Direction"
$epSheet.Range("B2").Value = "Exit"
$epSheet.Range("D1").Value = "Asset Name"
$epSheet.Range("D2").Value = ""

$epSheet.Columns.Item(2).ColumnWidth = 11.42578125
$epSheet.Columns.Item(4).ColumnWidth = 13.85546875

$epTable.Resize($epSheet.Range("A1:F2"))

# Fix up the table header/column names (resize keeps stale ones around for
# columns it thinks were already there).
$epTable.HeaderRowRange.Cells.Item(1, 1).Value = "ID"
$epTable.HeaderRowRange.Cells.Item(1, 2).Value = "Direction"
$epTable.HeaderRowRange.Cells.Item(1, 3).Value = "Description"
$epTable.HeaderRowRange.Cells.Item(1, 4).Value = "Asset Name"
$epTable.HeaderRowRange.Cells.Item(1, 5).Value = "Trust Level"
$epTable.HeaderRowRange.Cells.Item(1, 6).Value = "Microservice"

# -----------------------------------------------------------------------
# 2. Add the new "Vectors" worksheet at the end of the workbook with a
#    single-column table listing the known threat vectors.
# -----------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$vectorsSheet = $wb.Worksheets.Add($null, $lastSheet)
$vectorsSheet.Name = "Vectors"

$vectorsSheet.Range("A1").Value = "Name"
$vectorsSheet.Range("A2").Value = "Attack vector"

$vectorsTable = $vectorsSheet.ListObjects.Add(1, $vectorsSheet.Range("A1:A2"), $null, 1)
$vectorsTable.Name = "Table6"
$vectorsTable.TableStyle = "TableStyleMedium23"
